$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to match the latest scrape
$ws.Range('D2').Value = '43.981.47'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '2.304.80'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '113.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +17.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '270.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.624'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '48.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0953'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.06'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +14.36%  '
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.96'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('D15').Value = '2.643.69'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.858'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('D17').Value = '2.297.84'
$ws.Range('E17').Value = '  +0.32%  '
$ws.Range('D18').Value = '43.864.25'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('E20').Value = '  +9.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('E22').Value = '  -3.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.23'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.73'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.85'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.60%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.81%  '
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('E30').Value = '  +2.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.59'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0941'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.63'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.72'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.33%  '
$ws.Range('E35').Value = '  +1.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.69'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0364'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.29%  '
$ws.Range('E38').Value = '  +2.56%  '
$ws.Range('E39').Value = '  +9.00%  '
$ws.Range('E40').Value = '  +2.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.65'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +14.15%  '
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.46'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +24.01%  '
$ws.Range('B43').Value = 'LidoDAOToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.18%  '
$ws.Range('B44').Value = 'Celestia'
$ws.Range('C44').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.73'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +12.09%  '
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.40'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.64'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.21%  '
$ws.Range('E49').Value = '  -1.76%  '
$ws.Range('E50').Value = '  +3.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.467'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.93%  '
